# Updates the "Price" (D) and "Volume(1h)" (E) columns in the cryptos
# price table to the latest scrape, per the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.148.65"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "2.479.52"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D4").Value = "'1.00"  # leading ' keeps numeric-looking text as text
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'319.91"  # leading ' keeps numeric-looking text as text
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "'108.12"  # leading ' keeps numeric-looking text as text
$ws.Range("E6").Value = "  +3.06%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "'1.00"  # leading ' keeps numeric-looking text as text

$ws.Range("D9").Value = "'0.533"  # leading ' keeps numeric-looking text as text
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").Value = "'38.73"  # leading ' keeps numeric-looking text as text
$ws.Range("E10").Value = "  +7.28%  "

$ws.Range("D11").Value = "'0.0807"  # leading ' keeps numeric-looking text as text
$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "2.849.78"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "2.465.14"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "'0.844"  # leading ' keeps numeric-looking text as text
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "47.067.03"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").Value = "'12.67"  # leading ' keeps numeric-looking text as text
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("E20").Value = "  +2.01%  "

$ws.Range("D21").Value = "'2.75"  # leading ' keeps numeric-looking text as text
$ws.Range("E21").Value = "  +15.02%  "

$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("D23").Value = "'70.23"  # leading ' keeps numeric-looking text as text
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").Value = "'244.90"  # leading ' keeps numeric-looking text as text
$ws.Range("E24").Value = "  -1.84%  "

$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").Value = "'25.60"  # leading ' keeps numeric-looking text as text
$ws.Range("E27").Value = "  -2.11%  "

$ws.Range("E28").Value = "  +1.90%  "

$ws.Range("D29").Value = "'2.20"  # leading ' keeps numeric-looking text as text
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("E30").Value = "  +3.32%  "

$ws.Range("D31").Value = "'35.00"  # leading ' keeps numeric-looking text as text
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").Value = "'49.31"  # leading ' keeps numeric-looking text as text
$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("D33").Value = "'19.88"  # leading ' keeps numeric-looking text as text
$ws.Range("E33").Value = "  +1.30%  "

$ws.Range("D34").Value = "'5.33"  # leading ' keeps numeric-looking text as text
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("E36").Value = "  +0.26%  "

$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("D42").Value = "'118.59"  # leading ' keeps numeric-looking text as text
$ws.Range("E42").Value = "  -3.70%  "

$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").Value = "1.976.78"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("D47").Value = "'2.01"  # leading ' keeps numeric-looking text as text
$ws.Range("E47").Value = "  -4.36%  "

$ws.Range("D48").Value = "'9.04"  # leading ' keeps numeric-looking text as text
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("D50").Value = "'5.12"  # leading ' keeps numeric-looking text as text
$ws.Range("E50").Value = "  -4.64%  "

$ws.Range("D51").Value = "'57.00"  # leading ' keeps numeric-looking text as text
$ws.Range("E51").Value = "  +4.85%  "
